$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings used by B2 ("Hep 800M" -> "Deca 1500M") and C2 ("hej" -> "Linus")
$ws.Range("B2").Value = "Deca 1500M"
$ws.Range("C2").Value = "Linus"

# Update numeric values
$ws.Range("A2").Value = 4.0
$ws.Range("D2").Value = 400.0
$ws.Range("E2").Value = 124.0
